$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = -12.672
$ws.Range("A12").Value = -21.882
$ws.Range("C12").Value = -13.002
$ws.Range("C14").Value = -12.049
$ws.Range("C22").Value = -12.846
